$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 141, shifting existing row 141 (blank separator) and
# the summary rows (142-144) down to (142-145). This mirrors Excel's
# native "Insert Row" behavior, including shifting dependent formulas.
$ws.Rows(141).Insert()

# Populate the new data row 141 with the extra working-hours entry.
$ws.Range("A141").Value = 2014
$ws.Range("B141").Value = 7
$ws.Range("C141").Value = 14
$ws.Range("D141").Value = 0.33333333333333331
$ws.Range("E141").Value = 0.33333333333333331

# Extend the formula pattern in column F/G down through the new row 141
# (the existing shared formula covered F132:F140 / G132:G140; row 141 now
# continues the same per-row calculation).
$ws.Range("F141").Formula = "=(E141-D141)*24*60"
$ws.Range("G141").Formula = "=F141/60"

# The "sum [min]" total now needs to include the new row.
$ws.Range("F143").Formula = "=SUM(F2:F141)"

# Restore the active selection shown in the saved view.
$ws.Range("H143").Select()
